$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
WITH file_data AS (
    SELECT 
        file_name, 
        data_category, 
        file_description, 
        file_type, 
        file_size,
        file_access,  
        CAST("sample.id" AS TEXT) AS sample_id
    FROM df_sequencing_file
)
SELECT 
    fd.file_name AS "File Name",
    fd.data_category AS "Data Category",
    COALESCE(fd.file_description, '') AS "File Description",
    fd.file_type AS "File Type",
    CASE     
        WHEN fd.file_size >= 1024 * 1024 * 1024 THEN 
            CASE 
                WHEN ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
                THEN CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) || ' GB'
                ELSE ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
            END
        WHEN fd.file_size >= 1024 * 1024 THEN 
            CASE 
                WHEN ROUND(fd.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT) 
                THEN CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT) || ' MB'
                ELSE ROUND(fd.file_size / (1024.0 * 1024.0), 2) || ' MB'
            END
        WHEN fd.file_size >= 1024 THEN 
            CASE 
                WHEN ROUND(fd.file_size / 1024.0, 2) = CAST(ROUND(fd.file_size / 1024.0, 0) AS INT) 
                THEN CAST(ROUND(fd.file_size / 1024.0, 0) AS INT) || ' KB'
                ELSE ROUND(fd.file_size / 1024.0, 2) || ' KB'
            END
        ELSE 
            CASE 
                WHEN ROUND(fd.file_size, 2) = CAST(ROUND(fd.file_size, 0) AS INT) 
                THEN CAST(ROUND(fd.file_size, 0) AS INT) || ' Bytes'
                ELSE ROUND(fd.file_size, 2) || ' Bytes'
            END
    END AS "File Size",
    fd.file_access AS "File Access",
    std.dbgap_accession AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID"    
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON CAST(std.id AS TEXT) = CAST(prt."study.id" AS TEXT)
LEFT JOIN 
    df_sample smp ON CAST(prt.id AS TEXT) = CAST(smp."participant.id" AS TEXT)
LEFT JOIN
    df_diagnosis dgn ON CAST(prt.id AS TEXT) = CAST(dgn."participant.id" AS TEXT)
JOIN 
    file_data fd ON CAST(smp.id AS TEXT) = fd.sample_id
WHERE 
    std.dbgap_accession = 'phs002371' 
    AND prt.sex_at_birth = 'Male' 
    AND dgn.anatomic_site = 'C42.0 : Blood'
ORDER BY  
    fd.file_name
LIMIT 100;
'@

$ws.Range("B5").Value = $newQuery
$ws.Rows.Item(5).RowHeight = 409.5
